# Add the new daily ranking sheet "2025-07-28" after the last existing sheet
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-07-28"

# Header row
$ws.Cells.Item(1, 1).Value = "rank"
$ws.Cells.Item(1, 2).Value = "title"
$ws.Cells.Item(1, 3).Value = "author"
$ws.Cells.Item(1, 4).Value = "latest_episode"

# Ranking data: rank, title, author, latest_episode
$rankingData = @(
    @(1, '魔術師クノンは見えている', 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)', '第39話①'),
    @(2, 'ダンジョンの幼なじみ', '久真やすひさ(著者)', 'お休みイラスト'),
    @(3, 'モブ高生の俺でも冒険者になればリア充になれますか？', '原作：百均 漫画：さぎやまれん キャラクター原案：hai', '第30話'),
    @(4, 'クセ強彼女は床にいざなう', '須河篤志(著者)', '第14話前半'),
    @(5, '悪役令嬢、俺', '弥南せいら(漫画) 猫屋敷のあ(原作)', '第1話'),
    @(6, '異世界のんびり農家', '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)', '第302話'),
    @(7, '氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic', '漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca', '第8話'),
    @(8, 'ぽんドロイド！ はまさん', 'はれやまはれぞう(著者)', '第5話'),
    @(9, 'はずれスキル念動力（ただしレベルＭＡＸ）で無双する～手をかざすだけです。詠唱とか必殺技とかいりません。念じるだけで倒せます～', 'さとう うなぽっぽ トダフミト', '9話③'),
    @(10, '異世界おじさん', '殆ど死んでいる(著者)', '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'),
    @(11, '物語の黒幕に転生して', '瀬川はじめ(漫画) 結城涼(原作) なかむら(キャラクター原案)', '第32話'),
    @(12, 'リビルドワールド', '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)', '第71話③'),
    @(13, '転生貴族の異世界冒険録 ～カインのやりすぎギルド日記～', '原作：夜州 漫画：佐々木あかね・香本セトラ キャラクター原案：藻', '第48話'),
    @(14, 'サーシャちゃんとクラスメイトオタクくん', 'はぐはぐ(著者)', '第82話'),
    @(15, '幼女戦記', '東條チカ(漫画) カルロ・ゼン(原作) 篠月しのぶ(キャラクター原案)', '第百六章：ドードーバード航空戦Ⅰ'),
    @(16, '俺だけ使える古代魔法 ～基礎すら使えないと追放された俺の魔法は、実は1万年前に失われた伝説魔法でした～', '原作／アトハ 漫画／川上ニナ キャラクター原案／片倉響', '第15話'),
    @(17, 'ナナシの器用貧乏', '原作：高球 漫画：葉原辰月 キャラクター原案：KeG', '第5話'),
    @(18, '男子高校生だけどギャルにTSしました', '太陽まりい(著者)', '第18話後編'),
    @(19, '田舎の黒ギャルJKと結婚しました', 'カヅチ(著者)', '第17話後半'),
    @(20, '三枝さんはメガネ先輩と恋を描く', 'セレビィ量産型(著者)', '第21話前半'),
    @(21, 'モテはるねえ蘆屋くん', '栗原和明(原作) 寺井赤音(作画)', '第16話前半'),
    @(22, '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる', 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)', '「名もなき英雄譚」後半'),
    @(23, 'ゆーち〇ーばー', 'ともお', '第1話'),
    @(24, 'ワンパンマン', '原作/ＯＮＥ 作画/村田雄介', '206撃目'),
    @(25, '魔王令嬢の教育係 ～勇者学院を追放された平民教師は魔王の娘たちの家庭教師となる～', '原作：新人 漫画：松原剛 キャラクター原案：巻羊', '第29話'),
    @(26, '淫獄団地', '搾精研究所(原作) 丈山雄為(漫画)', '第49話（前編）'),
    @(27, 'ミルク搾りハンターの異世界搾乳記～農家の冴えない男があらゆる種族の地区Bを弄び虜にする～', '空詠大智(著者)', '第17話前半'),
    @(28, 'まんきつしたい常連さん', 'しんみりん(著者)', '第46話前編'),
    @(29, 'ソードアート・オンライン Re:Aincrad', '樹深(作画) 佐藤ミト(構成) 川原礫(原作) ａｂｅｃ(キャラクターデザイン)', '第40話②'),
    @(30, 'よくわからないけれど異世界に転生していたようです', '内々けやき あし カオミン', '第136話 よくわからないけれどスカウトされたみたいです（２）'),
    @(31, 'リアリスト魔王による聖域なき異世界改革', '鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)', '第67幕③'),
    @(32, '男嫌いな美人姉妹を名前も告げずに助けたら一体どうなる?', 'みょん(原作) 司馬淳子(漫画) ぎうにう(キャラクターデザイン)', '第23話'),
    @(33, 'ガルルガール', '原聡志(著者)', '第15話'),
    @(34, '最強ギフトで領地経営スローライフ ～辺境の村を開拓していたら英雄級の人材がわんさかやってきた！～', '原作：音速炒飯 漫画：眠田瞼 キャラクター原案：riritto', '第28話'),
    @(35, 'アラサーがVTuberになった話。', '犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)', '第24話'),
    @(36, '美人女上司滝沢さん', 'やんBARU(著者)', '第202話'),
    @(37, '転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く', '都尾琉(漫画) 御峰。(原作)', '第26話④'),
    @(38, '解雇された暗黒兵士(30代)のスローなセカンドライフ', '岡沢六十四 るれくちぇ sage・ジョー', '第71話(前編) ダリエルVS.滾り'),
    @(39, '天獄で悪魔がボクを魅惑する', '銀河味めてお(著者)', '第34話'),
    @(40, '老後に備えて異世界で８万枚の金貨を貯めます', 'FUNA 東西 モトエ恵介', '第120話　会談［その5］'),
    @(41, '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~', '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)', '第80話その3'),
    @(42, '王様ランキング200話～', '十日草輔（とおかそうすけ）', '第262話'),
    @(43, '愚かな天使は悪魔と踊る', 'アズマサワヨシ(著者)', '第100話③'),
    @(44, 'ひびワれこうじつ', 'ぺろぽねそす(著者)', '第1話'),
    @(45, 'ずっと好きだった幼馴染と付き合い始めたら一途ビッチの性欲ジャンキーだったんだがどうすりゃいいですか？', '原作：トラ子猫 作画：あらいぐま', '第3話'),
    @(46, '転生前はドラゴンでした', '甘井ヤドラキ', '第44話+お知らせ'),
    @(47, '両親の借金を肩代わりしてもらう条件は日本一可愛い女子高生と一緒に暮らすことでした。', '美月めいあ(漫画) 雨音恵(原作) ｋａｋａｏ(キャラクター原案)', '第36話'),
    @(48, '本物のカノジョにしたくなるまで、私で試していいよ。', 'wano(漫画) 有丈ほえる(原作) 緋月ひぐれ(キャラクター原案)', '第8話②'),
    @(49, 'わたし、二番目の彼女でいいから。', 'にの子(漫画) 西条陽(原作) Re岳(キャラクター原案)', '第10話③：カーテンのなか'),
    @(50, 'クラスで２番目に可愛い女の子と友だちになった', '尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)', '第34話①')
)

$r = 2
foreach ($row in $rankingData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Header formatting to match the other daily sheets (bold, centered/top, thin border)
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
